$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D writes to stay text (prices are formatted strings, e.g. "1.00", "0.116")
# by switching the cell to Text format before assigning the value, matching how the
# source data (inline strings) is stored. Column E values (e.g. "  -0.56%  ") already
# contain spaces/percent signs so Excel keeps them as text without this step.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '50.948.79'
$ws.Range('E2').Value = '  -0.56%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.945.55'
$ws.Range('E3').Value = '  -0.54%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '378.88'
$ws.Range('E5').Value = '  -1.08%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '101.08'
$ws.Range('E6').Value = '  -2.14%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.542'
$ws.Range('E7').Value = '  +0.23%  '
$ws.Range('E9').Value = '  -1.51%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.09'
$ws.Range('E10').Value = '  -1.53%  '
$ws.Range('E11').Value = '  -0.63%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0846'
$ws.Range('E12').Value = '  +0.60%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.407.91'
$ws.Range('E13').Value = '  -0.49%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '18.28'
$ws.Range('E14').Value = '  +1.23%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.75'
$ws.Range('E15').Value = '  +4.20%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '12.04'
$ws.Range('E16').Value = '  +68.65%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.947.12'
$ws.Range('E17').Value = '  +0.18%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.998'
$ws.Range('E18').Value = '  +0.78%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '50.933.52'
$ws.Range('E19').Value = '  -0.48%  '
$ws.Range('E20').Value = '  -4.68%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.38'
$ws.Range('E21').Value = '  -1.73%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0949'
$ws.Range('E22').Value = '  -0.57%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '69.40'
$ws.Range('E23').Value = '  +1.31%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '266.33'
$ws.Range('E24').Value = '  +1.49%  '
$ws.Range('E25').Value = '  +9.64%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.09'
$ws.Range('E26').Value = '  -3.86%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  +0.03%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.06'
$ws.Range('E28').Value = '  -9.60%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '25.58'
$ws.Range('E29').Value = '  -0.57%  '
$ws.Range('E30').Value = '  -4.37%  '
$ws.Range('E31').Value = '  -3.64%  '
$ws.Range('E32').Value = '  +2.73%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.05'
$ws.Range('E33').Value = '  +0.11%  '
$ws.Range('E34').Value = '  -0.03%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '33.42'
$ws.Range('E35').Value = '  -1.63%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0430'
$ws.Range('E36').Value = '  -6.20%  '
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.08'
$ws.Range('E38').Value = '  +2.95%  '
$ws.Range('E41').Value = '  +0.94%  '
$ws.Range('E42').Value = '  -1.45%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '119.57'
$ws.Range('E43').Value = '  -1.71%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.52'
$ws.Range('E44').Value = '  +8.46%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.28'
$ws.Range('E45').Value = '  -0.42%  '
$ws.Range('E46').Value = '  -1.50%  '
$ws.Range('E47').Value = '  -2.21%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.000.52'
$ws.Range('E48').Value = '  -0.67%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.260'
$ws.Range('E49').Value = '  -5.07%  '
$ws.Range('E50').Value = '  -10.23%  '
$ws.Range('E51').Value = '  +3.74%  '

# Rows 39/40 swap ranking order: Celestia <-> Stellar
$ws.Range('B39').Value = 'Stellar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.116'
$ws.Range('E39').Value = '  +0.37%  '

$ws.Range('B40').Value = 'Celestia'
$ws.Range('C40').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '16.58'
$ws.Range('E40').Value = '  -1.44%  '
